# Rename PersonListPanel -> TaskListPanel and PersonCard -> TaskCard on the
# single slide's diagram shapes, and refresh the stale "1/7/2017" date-field
# text (slide master, every custom layout, and the notes master) to
# "10/22/2018".

$p = $ppt.ActivePresentation

# --- Slide diagram shape renames ---------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $t = $shape.TextFrame.TextRange.Text
        if ($t -eq "PersonListPanel") {
            $shape.TextFrame.TextRange.Text = "TaskListPanel"
        } elseif ($t -eq "PersonCard") {
            $shape.TextFrame.TextRange.Text = "TaskCard"
        }
    }
}

# --- Helper: refresh a "1/7/2017" date placeholder's displayed text -----
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "1/7/2017") {
                $shp.TextFrame.TextRange.Text = "10/22/2018"
            }
        }
    }
}

# Slide master's own date placeholder
Update-DateShape $p.SlideMaster.Shapes

# Every custom (slide) layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

# Notes master: writing through NotesMaster.Shapes(...) on this host
# misroutes into the slide master's shape collection, so update the date
# field via the HeadersFooters object instead, which correctly targets the
# notes master part.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "10/22/2018"
